$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = "LoginTest001"
$ws.Range("B5").Value = "loginTest001"
$ws.Range("C5").Value = "loginPassword001"

$ws.Columns.Item("A:C").AutoFit()

$ws.Range("D5").Select()
